$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 459, shifting the existing rows (and all
# rows below them) down by two. This grows the used range from A1:R471
# to A1:R473, matching the new <dimension> in the target workbook.
$ws.Range("A459:A460").EntireRow.Insert()

# Populate the first new row (459) with the new weekly record.
$ws.Cells.Item(459, 1).Value = 9
$ws.Cells.Item(459, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(459, 3).Value = "Metropolitana"
$ws.Cells.Item(459, 4).Value = 44939
$ws.Cells.Item(459, 5).Value = 13
$ws.Cells.Item(459, 6).Value = 100112052
$ws.Cells.Item(459, 7).Value = "Albahaca"
$ws.Cells.Item(459, 8).Value = "Sin especificar"
$ws.Cells.Item(459, 9).Value = "Primera"
$ws.Cells.Item(459, 10).Value = 340
$ws.Cells.Item(459, 11).Value = 4000
$ws.Cells.Item(459, 12).Value = 5000
$ws.Cells.Item(459, 13).Value = 4500
$ws.Cells.Item(459, 14).Value = "`$/docena de matas"
$ws.Cells.Item(459, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(459, 16).Value = 750
$ws.Cells.Item(459, 17).Value = 6
$ws.Cells.Item(459, 18).Value = "Hortaliza"

# Populate the second new row (460) with the new weekly record.
$ws.Cells.Item(460, 1).Value = 9
$ws.Cells.Item(460, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(460, 3).Value = "Metropolitana"
$ws.Cells.Item(460, 4).Value = 44939
$ws.Cells.Item(460, 5).Value = 13
$ws.Cells.Item(460, 6).Value = 100112052
$ws.Cells.Item(460, 7).Value = "Albahaca"
$ws.Cells.Item(460, 8).Value = "Sin especificar"
$ws.Cells.Item(460, 9).Value = "Primera"
$ws.Cells.Item(460, 10).Value = 430
$ws.Cells.Item(460, 11).Value = 4000
$ws.Cells.Item(460, 12).Value = 5000
$ws.Cells.Item(460, 13).Value = 4500
$ws.Cells.Item(460, 14).Value = "`$/paquete"
$ws.Cells.Item(460, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(460, 16).Value = 750
$ws.Cells.Item(460, 17).Value = 6
$ws.Cells.Item(460, 18).Value = "Hortaliza"
